$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns (left-to-right) to make room for the "_4" metric columns
$ws.Columns("G:G").Insert()
$ws.Columns("K:K").Insert()
$ws.Columns("O:O").Insert()
$ws.Columns("S:S").Insert()

# Header row
# W1 falls outside the original A1:R5 used range, so copy formatting from a
# neighboring header cell before assigning its value (keeps the bold/border/
# centered style consistent with the rest of row 1).
$ws.Range("V1:V1").Copy($ws.Range("W1"))

$ws.Range("A1").Value = "model"
$ws.Range("B1").Value = "feature_selection_method"
$ws.Range("C1").Value = "best_params"
$ws.Range("D1").Value = "mse_1"
$ws.Range("E1").Value = "mse_2"
$ws.Range("F1").Value = "mse_3"
$ws.Range("G1").Value = "mse_4"
$ws.Range("H1").Value = "mae_1"
$ws.Range("I1").Value = "mae_2"
$ws.Range("J1").Value = "mae_3"
$ws.Range("K1").Value = "mae_4"
$ws.Range("L1").Value = "r2_1"
$ws.Range("M1").Value = "r2_2"
$ws.Range("N1").Value = "r2_3"
$ws.Range("O1").Value = "r2_4"
$ws.Range("P1").Value = "rmse_1"
$ws.Range("Q1").Value = "rmse_2"
$ws.Range("R1").Value = "rmse_3"
$ws.Range("S1").Value = "rmse_4"
$ws.Range("T1").Value = "explained_variance_score_1"
$ws.Range("U1").Value = "explained_variance_score_2"
$ws.Range("V1").Value = "explained_variance_score_3"
$ws.Range("W1").Value = "explained_variance_score_4"

# Existing rows 2-5: update best_params + numeric metric values (including the new *_4 columns)
# New rows 6-9: full new scenario rows ("*_corr" feature-selection variants)

# Row 2
$ws.Range("A2").Value = "random_forest"
$ws.Range("B2").Value = "rfe"
$ws.Range("C2").Value = "{'max_depth': 5, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Range("D2").Value = 40.04541689428554
$ws.Range("E2").Value = 214.6325086822381
$ws.Range("F2").Value = 119.8564870303817
$ws.Range("G2").Value = 8.793242478804412
$ws.Range("H2").Value = 5.136627284613351
$ws.Range("I2").Value = 11.42074789772591
$ws.Range("J2").Value = 6.414199043720438
$ws.Range("K2").Value = 2.270428510133486
$ws.Range("L2").Value = 0.4502982455585148
$ws.Range("M2").Value = 0.5557785214373216
$ws.Range("N2").Value = 0.3854610372806777
$ws.Range("O2").Value = 0.7620564316620487
$ws.Range("P2").Value = 6.328144822480404
$ws.Range("Q2").Value = 14.65034158926808
$ws.Range("R2").Value = 10.94789874954923
$ws.Range("S2").Value = 2.96534019613339
$ws.Range("T2").Value = 0.4504620211856515
$ws.Range("U2").Value = 0.5593789900154709
$ws.Range("V2").Value = 0.398797216794625
$ws.Range("W2").Value = 0.7632133834531998

# Row 3
$ws.Range("A3").Value = "random_forest"
$ws.Range("B3").Value = "mutual information"
$ws.Range("C3").Value = "{'max_depth': 5, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Range("D3").Value = 40.14218905947948
$ws.Range("E3").Value = 213.4599595244804
$ws.Range("F3").Value = 125.5264908354093
$ws.Range("G3").Value = 8.940881991106414
$ws.Range("H3").Value = 5.116789038759141
$ws.Range("I3").Value = 11.19745738608311
$ws.Range("J3").Value = 6.541019558356266
$ws.Range("K3").Value = 2.222201615572862
$ws.Range("L3").Value = 0.4489698581146124
$ws.Range("M3").Value = 0.5582053277200442
$ws.Range("N3").Value = 0.3563892836919627
$ws.Range("O3").Value = 0.7580613328723257
$ws.Range("P3").Value = 6.335786380511852
$ws.Range("Q3").Value = 14.6102689750901
$ws.Range("R3").Value = 11.20386053266504
$ws.Range("S3").Value = 2.990130764884107
$ws.Range("T3").Value = 0.4492700227126445
$ws.Range("U3").Value = 0.5616106088729953
$ws.Range("V3").Value = 0.3717605509266285
$ws.Range("W3").Value = 0.7599249758949355

# Row 4
$ws.Range("A4").Value = "random_forest"
$ws.Range("B4").Value = "random forest"
$ws.Range("C4").Value = "{'max_depth': 5, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Range("D4").Value = 41.17807625865857
$ws.Range("E4").Value = 203.7438493446589
$ws.Range("F4").Value = 126.7072263582291
$ws.Range("G4").Value = 9.07436488528409
$ws.Range("H4").Value = 5.189926296555458
$ws.Range("I4").Value = 11.07045318071319
$ws.Range("J4").Value = 6.588301788985013
$ws.Range("K4").Value = 2.263723825469939
$ws.Range("L4").Value = 0.4347502780738941
$ws.Range("M4").Value = 0.5783146059298432
$ws.Range("N4").Value = 0.3503353102991367
$ws.Range("O4").Value = 0.7544493096363837
$ws.Range("P4").Value = 6.417014590809232
$ws.Range("Q4").Value = 14.27388697393457
$ws.Range("R4").Value = 11.25643044478262
$ws.Range("S4").Value = 3.012368650295659
$ws.Range("T4").Value = 0.4348767080115156
$ws.Range("U4").Value = 0.5840649168823479
$ws.Range("V4").Value = 0.3677891057809866
$ws.Range("W4").Value = 0.7568776181649078

# Row 5
$ws.Range("A5").Value = "random_forest"
$ws.Range("B5").Value = "all"
$ws.Range("C5").Value = "{'max_depth': 5, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Range("D5").Value = 40.67302149982101
$ws.Range("E5").Value = 209.0643741551446
$ws.Range("F5").Value = 126.0678816918444
$ws.Range("G5").Value = 9.491560969275879
$ws.Range("H5").Value = 5.106866505615931
$ws.Range("I5").Value = 11.19928891612536
$ws.Range("J5").Value = 6.643890935514915
$ws.Range("K5").Value = 2.3379739349392
$ws.Range("L5").Value = 0.4416831435190194
$ws.Range("M5").Value = 0.5673028006233948
$ws.Range("N5").Value = 0.3536134157887506
$ws.Range("O5").Value = 0.7431600582412461
$ws.Range("P5").Value = 6.377540395781199
$ws.Range("Q5").Value = 14.4590585500974
$ws.Range("R5").Value = 11.22799544406055
$ws.Range("S5").Value = 3.080837705767034
$ws.Range("T5").Value = 0.4418635536516377
$ws.Range("U5").Value = 0.5730165550813762
$ws.Range("V5").Value = 0.3709540123000195
$ws.Range("W5").Value = 0.7453670890052116

# Row 6
$ws.Range("A6").Value = "random_forest"
$ws.Range("B6").Value = "rfe_corr"
$ws.Range("C6").Value = "{'max_depth': 5, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Range("D6").Value = 40.05655414927478
$ws.Range("E6").Value = 203.47848701386
$ws.Range("F6").Value = 123.5383768435975
$ws.Range("G6").Value = 9.771522327320266
$ws.Range("H6").Value = 5.094473333020188
$ws.Range("I6").Value = 10.95964807388055
$ws.Range("J6").Value = 6.518491389504591
$ws.Range("K6").Value = 2.356047213162978
$ws.Range("L6").Value = 0.4501453649274202
$ws.Range("M6").Value = 0.5788638221118001
$ws.Range("N6").Value = 0.366582920603631
$ws.Range("O6").Value = 0.7355843539785248
$ws.Range("P6").Value = 6.329024739189663
$ws.Range("Q6").Value = 14.26458856798401
$ws.Range("R6").Value = 11.11478190715398
$ws.Range("S6").Value = 3.125943429961628
$ws.Range("T6").Value = 0.4502747978767251
$ws.Range("U6").Value = 0.5824824333028951
$ws.Range("V6").Value = 0.3826327080480995
$ws.Range("W6").Value = 0.7373689439010525

# Row 7
$ws.Range("A7").Value = "random_forest"
$ws.Range("B7").Value = "mutual information_corr"
$ws.Range("C7").Value = "{'max_depth': 5, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Range("D7").Value = 40.2268175853126
$ws.Range("E7").Value = 200.7873083473967
$ws.Range("F7").Value = 125.6237779757079
$ws.Range("G7").Value = 8.681919101944013
$ws.Range("H7").Value = 5.09887530177093
$ws.Range("I7").Value = 10.76096389933431
$ws.Range("J7").Value = 6.451401761038793
$ws.Range("K7").Value = 2.228435533249391
$ws.Range("L7").Value = 0.4478081658977708
$ws.Range("M7").Value = 0.5844337116575751
$ws.Range("N7").Value = 0.3558904643141689
$ws.Range("O7").Value = 0.7650688223237919
$ws.Range("P7").Value = 6.342461476848921
$ws.Range("Q7").Value = 14.16994383712923
$ws.Range("R7").Value = 11.20820137112587
$ws.Range("S7").Value = 2.946509647352951
$ws.Range("T7").Value = 0.4481598540843266
$ws.Range("U7").Value = 0.5868783505891328
$ws.Range("V7").Value = 0.3709231820556166
$ws.Range("W7").Value = 0.7663632148446653

# Row 8
$ws.Range("A8").Value = "random_forest"
$ws.Range("B8").Value = "random forest_corr"
$ws.Range("C8").Value = "{'max_depth': 5, 'max_features': 'log2', 'n_estimators': 100}"
$ws.Range("D8").Value = 40.70295956872089
$ws.Range("E8").Value = 194.5403292514916
$ws.Range("F8").Value = 120.2025904088349
$ws.Range("G8").Value = 9.12013497668109
$ws.Range("H8").Value = 5.133752998733526
$ws.Range("I8").Value = 10.84182255368582
$ws.Range("J8").Value = 6.300779495866182
$ws.Range("K8").Value = 2.311065422163566
$ws.Range("L8").Value = 0.4412721849056456
$ws.Range("M8").Value = 0.597362984616135
$ws.Range("N8").Value = 0.383686464902844
$ws.Range("O8").Value = 0.753210779151593
$ws.Range("P8").Value = 6.379887112537406
$ws.Range("Q8").Value = 13.94777147975588
$ws.Range("R8").Value = 10.96369419533557
$ws.Range("S8").Value = 3.019956121648308
$ws.Range("T8").Value = 0.4415351882678781
$ws.Range("U8").Value = 0.6027332514360879
$ws.Range("V8").Value = 0.398366720820723
$ws.Range("W8").Value = 0.7558258233020325

# Row 9
$ws.Range("A9").Value = "random_forest"
$ws.Range("B9").Value = "all_corr"
$ws.Range("C9").Value = "{'max_depth': 5, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Range("D9").Value = 40.12007187002515
$ws.Range("E9").Value = 211.3727720449667
$ws.Range("F9").Value = 125.116086554535
$ws.Range("G9").Value = 9.866569880249202
$ws.Range("H9").Value = 5.083135340004559
$ws.Range("I9").Value = 11.41540945082701
$ws.Range("J9").Value = 6.579528485003495
$ws.Range("K9").Value = 2.391660088193357
$ws.Range("L9").Value = 0.4492734598445794
$ws.Range("M9").Value = 0.5625251463434184
$ws.Range("N9").Value = 0.3584935454412672
$ws.Range("O9").Value = 0.7330123842005715
$ws.Range("P9").Value = 6.334040722163471
$ws.Range("Q9").Value = 14.53866472702933
$ws.Range("R9").Value = 11.18553023126463
$ws.Range("S9").Value = 3.141109657469666
$ws.Range("T9").Value = 0.4493366191766381
$ws.Range("U9").Value = 0.5668259992564464
$ws.Range("V9").Value = 0.3732732254189337
$ws.Range("W9").Value = 0.7351386177576507
